$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting existing rows 98:113 down to 99:114
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with a new weekly price record (same constant
# metadata as the surrounding rows, new date + price figures)
$ws.Range("A98").Value = 9
$ws.Range("B98").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 45154
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = 100112029
$ws.Range("G98").Value = "Orégano"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 16
$ws.Range("K98").Value = 21000
$ws.Range("L98").Value = 21000
$ws.Range("M98").Value = 21000
$ws.Range("N98").Value = "$/docena de atados"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 7000
$ws.Range("Q98").Value = 3
$ws.Range("R98").Value = "Hortaliza"

# Match the date-cell style used by the rest of column D
$ws.Range("D98").NumberFormat = $ws.Range("D99").NumberFormat
